$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.536.71'

$ws.Range("E2").Value = '  +2.04%  '

$ws.Range("D3").Value = '1.564.88'

$ws.Range("E3").Value = '  +0.36%  '

$ws.Range("E4").Value = '  -1.57%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.54'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +1.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.489'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  -0.13%  '

$ws.Range("E7").Value = '  -1.58%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.48'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  +2.04%  '

$ws.Range("E9").Value = '  +0.58%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0595'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  -0.26%  '

$ws.Range("D12").Value = '1.788.94'

$ws.Range("E12").Value = '  +0.47%  '

$ws.Range("D13").Value = '1.576.14'

$ws.Range("E13").Value = '  +4.81%  '

$ws.Range("E14").Value = '  +0.72%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.521'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +0.30%  '

$ws.Range("D16").Value = '27.505.97'

$ws.Range("E16").Value = '  +1.94%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.25'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +0.75%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '224.54'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  +4.03%  '

$ws.Range("E19").Value = '  +2.04%  '

$ws.Range("E20").Value = '  -0.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.989'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -1.57%  '

$ws.Range("E22").Value = '  +0.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.38'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +1.91%  '

$ws.Range("E24").Value = '  +0.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.47'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  -1.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.18'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  +1.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.61'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  +0.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.108'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  +2.03%  '

$ws.Range("E29").Value = '  -1.60%  '

$ws.Range("E30").Value = '  +1.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0471'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  -0.64%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.24'
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").Value = '1.465.71'

$ws.Range("E33").Value = '  +3.00%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.19'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +0.34%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.11'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +3.30%  '

$ws.Range("E36").Value = '  +1.21%  '

$ws.Range("E37").Value = '  -1.11%  '

$ws.Range("E38").Value = '  +0.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.542'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +1.83%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.817'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  +0.99%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.91'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +9.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.71'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  -1.54%  '

$ws.Range("E43").Value = '  +1.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.989'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  -1.59%  '

$ws.Range("E45").Value = '  -2.80%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.22'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  +0.94%  '

$ws.Range("D47").Value = '1.703.29'

$ws.Range("E47").Value = '  +0.70%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.50'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  -0.26%  '

$ws.Range("E49").Value = '  +1.07%  '

$ws.Range("E50").Value = '  -2.17%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0950'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  -0.92%  '

Write-Host "Applied 80 changes"